# Refresh the coin price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.903.59"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.633.83"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Value = "'214.42"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "'0.5082"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'0.2553"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'19.48"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "'0.07761"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "'4.280"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "1.647.59"
$ws.Range("E13").Value = "  -2.36%  "
$ws.Range("D14").Value = "'0.5424"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0₅7719"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'64.08"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "25.911.61"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'195.69"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'4.417"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("D21").Value = "'9.902"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "'6.020"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "'1.005"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'1.868"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").Value = "'141.13"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").Value = "'0.1192"
$ws.Range("E26").Value = "  +5.88%  "
$ws.Range("D27").Value = "'6.818"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").Value = "'15.60"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "'0.04918"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "'3.237"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Value = "'3.169"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").Value = "'1.520"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "'2.366"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").Value = "'0.8892"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").Value = "'2.576"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "1.136.09"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("D38").Value = "'0.5404"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").Value = "'0.01548"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").Value = "0.0₈128"
$ws.Range("E42").Value = "  +7.40%  "
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'98.64"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.432"
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("D46").Value = "1.770.53"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "'0.4519"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'0.9972"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "'54.69"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").Value = "'0.05045"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  -0.40%  "
